$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "third_trait" row (row 4) content by clearing rows first,
# then rewrite the table with the new custom-traits header/content layout.
$ws.Cells.Clear()

# Row 1: header row, bold
$ws.Range("A1").Value = "trait"
$ws.Range("B1").Value = "formula"
$ws.Range("A1:B1").Font.Bold = $true

# Row 2
$ws.Range("A2").Value = "IgGI_first_trait"
$ws.Range("B2").Value = "0.5 * (IgGI1H4N4F1 + IgGI1H5N4F1)"

# Row 3
$ws.Range("A3").Value = "second_trait"
$ws.Range("B3").Value = "IgGI1H4N4F1S1 / (IgGI1H4N4F1 + IgGI1H4N5F1S1)"

# Column widths to match target layout (engine snaps ColumnWidth to 1/6-character
# increments, so pick inputs landing on the grid points nearest the target widths)
$ws.Columns.Item(1).ColumnWidth = 14.833333333333334
$ws.Columns.Item(2).ColumnWidth = 41.666666666666664

# Remove the selection anchor discrepancy by selecting A1
$ws.Range("A1").Select()

# Page setup so a pageSetup element is emitted
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
